# Fruta / hortaliza, semanal
# Insert two new weekly observation rows at the top of the data block
# (row 36 onward), pushing the existing rows down by two. Excel's
# Rows(...).Insert() shifts everything below automatically, which is
# exactly what the target workbook shows (old row 36 -> new row 38,
# ..., old row 68 -> new row 70), and it grows the used range from
# A1:T68 to A1:T70.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("36:37").Insert()

# --- New row 36: Black Amber / Primera, week of 2023-01-18 (44944) ---
$ws.Cells.Item(36, 1).Value = 7
$ws.Cells.Item(36, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(36, 3).Value = "Ñuble"
$ws.Cells.Item(36, 4).Value = 44944
$ws.Cells.Item(36, 5).Value = 16
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100103
$ws.Cells.Item(36, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(36, 9).Value = 100103002
$ws.Cells.Item(36, 10).Value = "Ciruela"
$ws.Cells.Item(36, 11).Value = "Black Amber"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 60
$ws.Cells.Item(36, 14).Value = 14000
$ws.Cells.Item(36, 15).Value = 15000
$ws.Cells.Item(36, 16).Value = 14500
$ws.Cells.Item(36, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(36, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(36, 19).Value = 806
$ws.Cells.Item(36, 20).Value = 18

# --- New row 37: Black Amber / Segunda, week of 2023-01-18 (44944) ---
$ws.Cells.Item(37, 1).Value = 7
$ws.Cells.Item(37, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(37, 3).Value = "Ñuble"
$ws.Cells.Item(37, 4).Value = 44944
$ws.Cells.Item(37, 5).Value = 16
$ws.Cells.Item(37, 6).Value = "Fruta"
$ws.Cells.Item(37, 7).Value = 100103
$ws.Cells.Item(37, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(37, 9).Value = 100103002
$ws.Cells.Item(37, 10).Value = "Ciruela"
$ws.Cells.Item(37, 11).Value = "Black Amber"
$ws.Cells.Item(37, 12).Value = "Segunda"
$ws.Cells.Item(37, 13).Value = 80
$ws.Cells.Item(37, 14).Value = 11000
$ws.Cells.Item(37, 15).Value = 12000
$ws.Cells.Item(37, 16).Value = 11500
$ws.Cells.Item(37, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(37, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(37, 19).Value = 639
$ws.Cells.Item(37, 20).Value = 18
